# Update forecast-error table: refresh evaluation metrics for rows Q0..Q8
# (B2:G10) and append a new evaluation row Q9 (row 11).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated values for existing rows (row -> B,C,D,E,F,G) ---
$rows = @{
    2  = @(0.2134401944279581, 1.107048576932442,  5.469072502632633, 2.338604819680451, 2.352017453539809, 51)
    3  = @(0.2381905613280194, 1.157832168035669,  5.127790050230712, 2.264462419699367, 2.27476289563374,  50)
    4  = @(0.2567033173151257, 1.116656496429176,  4.777664275584832, 2.185786877896569, 2.193155149364095, 49)
    5  = @(0.2256811907542642, 1.135244597037652,  5.196387630754981, 2.279558648237632, 2.292364153068935, 48)
    6  = @(0.2532833917168094, 1.231919099044784,  5.524199383542229, 2.350361543155059, 2.361936404502555, 47)
    7  = @(0.1912077111841825, 1.218487072762304,  5.000023379858966, 2.236073205389074, 2.25250131378905,  46)
    8  = @(0.1745139725577196, 1.207597693525557,  5.537112892574842, 2.353107072059162, 2.373143290831724, 45)
    9  = @(0.2127461743163057, 1.208769835132661,  5.244261527796105, 2.290035267806176, 2.306492485061966, 44)
    10 = @(0.2033757677206456, 1.265673175597562,  5.739404322570661, 2.395705391439161, 2.415307496178391, 43)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
    $ws.Cells.Item($r, 5).Value = $vals[3]
    $ws.Cells.Item($r, 6).Value = $vals[4]
    $ws.Cells.Item($r, 7).Value = $vals[5]
}

# --- Append the new "Q9" row (row 11), matching the formatting of row 10 ---
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A11").Value = "Q9"

$ws.Cells.Item(11, 2).Value = 0.1762107366466173
$ws.Cells.Item(11, 3).Value = 1.24023323203181
$ws.Cells.Item(11, 4).Value = 5.544438103667789
$ws.Cells.Item(11, 5).Value = 2.354663055230576
$ws.Cells.Item(11, 6).Value = 2.376522828646756
$ws.Cells.Item(11, 7).Value = 42
